# Updates loading_percent.xlsx values for Case_3_146 (380 kV case).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.88798970853751
$ws.Range("C2").Value = 6.574817944023444
$ws.Range("D2").Value = 7.665164441200761
$ws.Range("F2").Value = 44.54936210137549
$ws.Range("G2").Value = 3.643946163734529
$ws.Range("M2").Value = 19.97276543822908

$ws.Range("B3").Value = 12.55740894293622
$ws.Range("C3").Value = 6.126145337998938
$ws.Range("D3").Value = 7.48964434358138
$ws.Range("F3").Value = 43.05163983448085
$ws.Range("G3").Value = 3.64865037905393
$ws.Range("M3").Value = 19.37634560475837

$ws.Range("B4").Value = 12.35774551085176
$ws.Range("C4").Value = 5.833934163621654
$ws.Range("D4").Value = 7.379735629858411
$ws.Range("F4").Value = 42.11140927492345
$ws.Range("G4").Value = 3.651679232675459
$ws.Range("M4").Value = 19.00932864804849

$ws.Range("B5").Value = 12.27738491660901
$ws.Range("C5").Value = 5.710634602709279
$ws.Range("D5").Value = 7.334445440941832
$ws.Range("F5").Value = 41.72357098524044
$ws.Range("G5").Value = 3.652949006992141
$ws.Range("M5").Value = 18.8598085995581

$ws.Range("B6").Value = 12.26410675540896
$ws.Range("C6").Value = 5.689905401646737
$ws.Range("D6").Value = 7.326895842909241
$ws.Range("F6").Value = 41.65890307469879
$ws.Range("G6").Value = 3.653162000801406
$ws.Range("M6").Value = 18.83499089514115

$ws.Range("B7").Value = 12.35665744819888
$ws.Range("C7").Value = 5.832288414657173
$ws.Range("D7").Value = 7.379126811938118
$ws.Range("F7").Value = 42.10619703808617
$ws.Range("G7").Value = 3.65169621335744
$ws.Range("M7").Value = 19.00731166260864

$ws.Range("B8").Value = 12.77341889958291
$ws.Range("C8").Value = 6.423577308784815
$ws.Range("D8").Value = 7.605110233849905
$ws.Range("F8").Value = 44.03751416307942
$ws.Range("G8").Value = 3.645539133603154
$ws.Range("M8").Value = 19.76743290242693

$ws.Range("B9").Value = 13.60968510619682
$ws.Range("C9").Value = 7.451104347692016
$ws.Range("D9").Value = 8.02979216273217
$ws.Range("F9").Value = 47.64142131292402
$ws.Range("G9").Value = 3.634571459225055
$ws.Range("M9").Value = 21.24139251680351

$ws.Range("B10").Value = 14.22632265512459
$ws.Range("C10").Value = 8.126554599458533
$ws.Range("D10").Value = 8.328675744415522
$ws.Range("F10").Value = 50.15394026815972
$ws.Range("G10").Value = 3.627176658636887
$ws.Range("M10").Value = 22.30158717517285

$ws.Range("B11").Value = 14.50551645540473
$ws.Range("C11").Value = 8.416789611111769
$ws.Range("D11").Value = 8.461439683325246
$ws.Range("F11").Value = 51.26343786854838
$ws.Range("G11").Value = 3.623954126003767
$ws.Range("M11").Value = 22.77662997845971

$ws.Range("B12").Value = 14.6109127140456
$ws.Range("C12").Value = 8.524261740093596
$ws.Range("D12").Value = 8.511226864601994
$ws.Range("F12").Value = 51.67847008162131
$ws.Range("G12").Value = 3.622753978087038
$ws.Range("M12").Value = 22.95530315406542

$ws.Range("B13").Value = 14.58823044666816
$ws.Range("C13").Value = 8.501223682568821
$ws.Range("D13").Value = 8.500526407847312
$ws.Range("F13").Value = 51.58931713957915
$ws.Range("G13").Value = 3.62301155784223
$ws.Range("M13").Value = 22.91687931497598

$ws.Range("B14").Value = 14.51419480786614
$ws.Range("C14").Value = 8.42568009196092
$ws.Range("D14").Value = 8.465545627782712
$ws.Range("F14").Value = 51.29768682526353
$ws.Range("G14").Value = 3.62385498619292
$ws.Range("M14").Value = 22.79135481524952

$ws.Range("B15").Value = 14.46879913668568
$ws.Range("C15").Value = 8.37909099604213
$ws.Range("D15").Value = 8.44405460435148
$ws.Range("F15").Value = 51.11838104288038
$ws.Range("G15").Value = 3.624374229927586
$ws.Range("M15").Value = 22.71430432356053

$ws.Range("B16").Value = 14.20803805385397
$ws.Range("C16").Value = 8.107245676529759
$ws.Range("D16").Value = 8.319932520480853
$ws.Range("F16").Value = 50.08073231401126
$ws.Range("G16").Value = 3.62739008885858
$ws.Range("M16").Value = 22.27038100697327

$ws.Range("B17").Value = 14.04762991580623
$ws.Range("C17").Value = 7.936126387933118
$ws.Range("D17").Value = 8.242947855206831
$ws.Range("F17").Value = 49.43537432401035
$ws.Range("G17").Value = 3.629276310131734
$ws.Range("M17").Value = 21.99606396260711

$ws.Range("B18").Value = 13.95525101846533
$ws.Range("C18").Value = 7.83609916788793
$ws.Range("D18").Value = 8.198368650565719
$ws.Range("F18").Value = 49.06104757204588
$ws.Range("G18").Value = 3.630374534570278
$ws.Range("M18").Value = 21.83761138271729

$ws.Range("B19").Value = 13.92395752773659
$ws.Range("C19").Value = 7.801955723110201
$ws.Range("D19").Value = 8.183224304740879
$ws.Range("F19").Value = 48.93377841621393
$ws.Range("G19").Value = 3.630748667591777
$ws.Range("M19").Value = 21.78385209508432

$ws.Range("B20").Value = 14.06471872050128
$ws.Range("C20").Value = 7.95450821728539
$ws.Range("D20").Value = 8.25117422426332
$ws.Range("F20").Value = 49.50440041064606
$ws.Range("G20").Value = 3.629074141389045
$ws.Range("M20").Value = 22.02533640981761

$ws.Range("B21").Value = 14.53595083573144
$ws.Range("C21").Value = 8.447934999337338
$ws.Range("D21").Value = 8.475833774034388
$ws.Range("F21").Value = 51.38348661984872
$ws.Range("G21").Value = 3.623606705254728
$ws.Range("M21").Value = 22.82825868127059

$ws.Range("B22").Value = 14.84195738781074
$ws.Range("C22").Value = 8.756240352644914
$ws.Range("D22").Value = 8.61980817324199
$ws.Range("F22").Value = 52.58167642993362
$ws.Range("G22").Value = 3.620150816761483
$ws.Range("M22").Value = 23.34586313992361

$ws.Range("B23").Value = 14.67885939075814
$ws.Range("C23").Value = 8.592984142314112
$ws.Range("D23").Value = 8.543235972640241
$ws.Range("F23").Value = 51.94500565447632
$ws.Range("G23").Value = 3.621984606619376
$ws.Range("M23").Value = 23.07031514640595

$ws.Range("B24").Value = 14.05699335375713
$ws.Range("C24").Value = 7.946202927993022
$ws.Range("D24").Value = 8.247456076479532
$ws.Range("F24").Value = 49.47320397832867
$ws.Range("G24").Value = 3.629165498847502
$ws.Range("M24").Value = 22.01210464065968

$ws.Range("B25").Value = 13.38241367682543
$ws.Range("C25").Value = 7.187159102953592
$ws.Range("D25").Value = 7.917082684610495
$ws.Range("F25").Value = 46.6886871282603
$ws.Range("G25").Value = 3.637421230669551
$ws.Range("M25").Value = 20.84578926037821
